# Apply updated cryptocurrency price/volume data to Sheet1 (columns D and E)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.140.28"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.902.90"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.88"
$ws.Range("E5").Value = "  -0.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9995"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5218"
$ws.Range("E7").Value = "  +1.34%  "
$ws.Range("E8").Value = "  +0.40%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07250"
$ws.Range("E9").Value = "  +0.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.13"
$ws.Range("E10").Value = "  -0.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9019"
$ws.Range("E11").Value = "  -0.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08512"
$ws.Range("E12").Value = "  +11.43%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.900.45"
$ws.Range("E13").Value = "  +0.43%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "95.04"
$ws.Range("E14").Value = "  +0.32%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.288"
$ws.Range("E15").Value = "  +0.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.000"
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("E17").Value = "  +1.44%  "
$ws.Range("E18").Value = "  +0.42%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9999"
$ws.Range("E19").Value = "  +0.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.172.92"
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("E21").Value = "  -0.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.155.44"
$ws.Range("E22").Value = "  +1.58%  "
$ws.Range("E23").Value = "  +0.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.423"
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.96"
$ws.Range("E25").Value = "  +0.41%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.284"
$ws.Range("E26").Value = "  +3.50%  "
$ws.Range("E27").Value = "  -2.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.20"
$ws.Range("E28").Value = "  +0.67%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.92"
$ws.Range("E29").Value = "  +0.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.806"
$ws.Range("E30").Value = "  -1.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.900"
$ws.Range("E31").Value = "  -1.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09251"
$ws.Range("E32").Value = "  +0.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.8019"
$ws.Range("E33").Value = "  +4.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05051"
$ws.Range("E34").Value = "  -0.66%  "
$ws.Range("E35").Value = "  -0.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.450"
$ws.Range("E36").Value = "  +4.65%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.939"
$ws.Range("E37").Value = "  -1.20%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.621"
$ws.Range("E38").Value = "  +0.35%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5706"
$ws.Range("E39").Value = "  +1.59%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01994"
$ws.Range("E40").Value = "  -0.37%  "
$ws.Range("E41").Value = "  -0.22%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.042"
$ws.Range("E42").Value = "  +0.85%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.630"
$ws.Range("E43").Value = "  -0.49%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "116.29"
$ws.Range("E44").Value = "  -1.27%  "
$ws.Range("E45").Value = "  +0.13%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4862"
$ws.Range("E46").Value = "  +0.98%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9991"
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("E48").Value = "  -1.17%  "
$ws.Range("E49").Value = "  +1.45%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.47"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.96"
$ws.Range("E51").Value = "  +0.03%  "
